$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wedge A data: move raw MCNP values from column S (A-Sab) to column R (A)
# Row 10: only clear the stray S10 value (R10 already holds its value).
$ws.Range("S10").ClearContents()

# Row 11
$ws.Range("R11").Value = 0.0136516
$ws.Range("S11").ClearContents()

# Row 12
$ws.Range("R12").Value = 0.040834
$ws.Range("S12").ClearContents()

# Row 13
$ws.Range("R13").Value = 0.389647
$ws.Range("S13").ClearContents()

# Row 14
$ws.Range("R14").Value = 0.744999
$ws.Range("S14").ClearContents()

# Rows 15-22: just clear the stale S values (no corresponding R value given yet)
$ws.Range("S15").ClearContents()
$ws.Range("S16").ClearContents()
$ws.Range("S17").ClearContents()
$ws.Range("S18").ClearContents()
$ws.Range("S19").ClearContents()
$ws.Range("S20").ClearContents()
$ws.Range("S21").ClearContents()
$ws.Range("S22").ClearContents()

# --- View state: scroll sheet back to top and move the active selection
$ws.Activate()
$ws.Range("L5").Select()
